# Add a new quiz row (row 29) to the "Csillagászat" sheet, matching the
# previously-blank template row's structure used by the rows above it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Value = "Melyik éghajlat jellemző a Dél-Kínai-hegyvidékre?"
$ws.Range("B29").Value = "tobbszoros"
$ws.Range("C29").Value = "mérsékelt övezeti monszun"
$ws.Range("D29").Value = "forró övezeti monszun;szavanna;egyenlítői"
$ws.Range("F29").Value = 4
$ws.Range("G29").Value = "közepes"
